# Applies the "Designed some new stages" edit:
#  1. Notes Master date placeholder field text: 3/14/2019 -> 7/30/2019
#  2. Stage table, "START POSITION" cell: trim the end of the existing
#     sentence and append the new "wrists below belt" wording.
#  3. Stage table, "GUN READY CONDITION" cell: PCC is now also loaded.

$p = $ppt.ActivePresentation

# --- 1. Notes master date field -------------------------------------------
$nm = $p.NotesMaster
$dateAndTime = $nm.HeadersFooters.DateAndTime
$dateAndTime.Text = "7/30/2019"

# --- locate the stage-info table on slide 1 --------------------------------
$s = $p.Slides.Item(1)
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}
$tbl = $tableShape.Table

# --- 2. START POSITION cell -------------------------------------------------
$startCell = $tbl.Cell(2, 1)
$startRange = $startCell.Shape.TextFrame.TextRange
$startText = $startRange.Text
$oldTail = "heels against rear fault line"
$newTail = "wrists below belt"
$idx = $startText.IndexOf($oldTail)
if ($idx -ge 0) {
    $sub = $startRange.Characters($idx + 1, $oldTail.Length)
    $sub.Text = $newTail
}

# --- 3. GUN READY CONDITION cell -------------------------------------------
$gunCell = $tbl.Cell(3, 1)
$gunRange = $gunCell.Shape.TextFrame.TextRange
$gunText = $gunRange.Text
$oldPhrase = "PCC shouldered"
$newPhrase = "PCC loaded, shouldered"
$idx2 = $gunText.IndexOf($oldPhrase)
if ($idx2 -ge 0) {
    $sub2 = $gunRange.Characters($idx2 + 1, $oldPhrase.Length)
    $sub2.Text = $newPhrase
}
